$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44539
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 3800
$ws.Range("O2").Value = 4000
$ws.Range("P2").Value = 3900
$ws.Range("R2").Value = "Región del Maule"
$ws.Range("S2").Value = 1950
$ws.Range("D3").Value = 44974
$ws.Range("M3").Value = 40
$ws.Range("N3").Value = 3000
$ws.Range("O3").Value = 3000
$ws.Range("P3").Value = 3000
$ws.Range("S3").Value = 1500
$ws.Range("D4").Value = 44974
$ws.Range("D5").Value = 44965
$ws.Range("M5").Value = 50
$ws.Range("N5").Value = 3000
$ws.Range("P5").Value = 3000
$ws.Range("S5").Value = 1500
$ws.Range("D6").Value = 44174
$ws.Range("M6").Value = 150
$ws.Range("N6").Value = 3700
$ws.Range("O6").Value = 3800
$ws.Range("P6").Value = 3747
$ws.Range("R6").Value = "Provincia de Linares"
$ws.Range("S6").Value = 1874
$ws.Range("D7").Value = 44944
$ws.Range("M7").Value = 60
$ws.Range("N7").Value = 2500
$ws.Range("O7").Value = 2500
$ws.Range("P7").Value = 2500
$ws.Range("S7").Value = 1250
$ws.Range("D8").Value = 45006
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 60
$ws.Range("N8").Value = 4000
$ws.Range("O8").Value = 4000
$ws.Range("P8").Value = 4000
$ws.Range("S8").Value = 2000
$ws.Range("D9").Value = 44540
$ws.Range("M9").Value = 240
$ws.Range("N9").Value = 3500
$ws.Range("O9").Value = 3800
$ws.Range("P9").Value = 3650
$ws.Range("R9").Value = "Región del Maule"
$ws.Range("S9").Value = 1825
$ws.Range("D10").Value = 44960
$ws.Range("L10").Value = "Segunda"
$ws.Range("D11").Value = 44967
$ws.Range("M11").Value = 50
$ws.Range("N11").Value = 3000
$ws.Range("P11").Value = 3000
$ws.Range("R11").Value = "Provincia de Diguillín"
$ws.Range("S11").Value = 1500
$ws.Range("D12").Value = 44967
$ws.Range("L12").Value = "Segunda"
$ws.Range("M12").Value = 30
$ws.Range("N12").Value = 2500
$ws.Range("O12").Value = 2500
$ws.Range("P12").Value = 2500
$ws.Range("Q12").Value = "`$/bandeja 2 kilos"
$ws.Range("S12").Value = 1250
$ws.Range("T12").Value = 2
$ws.Range("D13").Value = 44985
$ws.Range("M13").Value = 50
$ws.Range("D14").Value = 44985
$ws.Range("M14").Value = 50
$ws.Range("D15").Value = 45009
$ws.Range("N15").Value = 4000
$ws.Range("O15").Value = 4000
$ws.Range("P15").Value = 4000
$ws.Range("S15").Value = 2000
$ws.Range("D16").Value = 44988
$ws.Range("M16").Value = 30
$ws.Range("D17").Value = 44988
$ws.Range("L17").Value = "Segunda"
$ws.Range("M17").Value = 30
$ws.Range("D18").Value = 44994
$ws.Range("M18").Value = 60
$ws.Range("N18").Value = 3000
$ws.Range("O18").Value = 3200
$ws.Range("P18").Value = 3100
$ws.Range("R18").Value = "Provincia de Diguillín"
$ws.Range("S18").Value = 1550
$ws.Range("D19").Value = 44181
$ws.Range("M19").Value = 65
$ws.Range("N19").Value = 3600
$ws.Range("O19").Value = 3800
$ws.Range("P19").Value = 3692
$ws.Range("S19").Value = 1846
$ws.Range("D20").Value = 44181
$ws.Range("M20").Value = 80
$ws.Range("N20").Value = 1800
$ws.Range("O20").Value = 2000
$ws.Range("P20").Value = 1875
$ws.Range("Q20").Value = "`$/envase 1 kilo"
$ws.Range("R20").Value = "Provincia de Diguillín"
$ws.Range("S20").Value = 1875
$ws.Range("T20").Value = 1
$ws.Range("D21").Value = 44971
$ws.Range("D22").Value = 44992
$ws.Range("L22").Value = "Primera"
$ws.Range("M22").Value = 100
$ws.Range("N22").Value = 3000
$ws.Range("O22").Value = 3000
$ws.Range("P22").Value = 3000
$ws.Range("S22").Value = 1500
$ws.Range("D23").Value = 44596
$ws.Range("M23").Value = 120
$ws.Range("N23").Value = 2500
$ws.Range("O23").Value = 2700
$ws.Range("P23").Value = 2600
$ws.Range("R23").Value = "Provincia de Linares"
$ws.Range("S23").Value = 1300
$ws.Range("D24").Value = 44956
$ws.Range("M24").Value = 50
$ws.Range("N24").Value = 3000
$ws.Range("O24").Value = 3000
$ws.Range("P24").Value = 3000
$ws.Range("R24").Value = "Provincia de Diguillín"
$ws.Range("S24").Value = 1500
$ws.Range("D25").Value = 44970
$ws.Range("M25").Value = 50
$ws.Range("O25").Value = 3000
$ws.Range("P25").Value = 3000
$ws.Range("S25").Value = 1500
$ws.Range("D26").Value = 44970
$ws.Range("L26").Value = "Segunda"
$ws.Range("M26").Value = 30
$ws.Range("N26").Value = 2500
$ws.Range("O26").Value = 2500
$ws.Range("P26").Value = 2500
$ws.Range("S26").Value = 1250
$ws.Range("D27").Value = 44963
$ws.Range("L27").Value = "Primera"
$ws.Range("N27").Value = 3000
$ws.Range("O27").Value = 3000
$ws.Range("P27").Value = 3000
$ws.Range("S27").Value = 1500
$ws.Range("D28").Value = 44963
$ws.Range("L28").Value = "Segunda"
$ws.Range("N28").Value = 2500
$ws.Range("O28").Value = 2500
$ws.Range("P28").Value = 2500
$ws.Range("S28").Value = 1250
$ws.Range("D29").Value = 44949
$ws.Range("M29").Value = 60
$ws.Range("N29").Value = 2800
$ws.Range("P29").Value = 2900
$ws.Range("S29").Value = 1450
$ws.Range("D30").Value = 45008
$ws.Range("D31").Value = 44966
$ws.Range("L31").Value = "Segunda"
$ws.Range("M31").Value = 30
$ws.Range("N31").Value = 2500
$ws.Range("O31").Value = 2500
$ws.Range("P31").Value = 2500
$ws.Range("S31").Value = 1250
$ws.Range("D32").Value = 44931
$ws.Range("M32").Value = 100
$ws.Range("D33").Value = 44942
$ws.Range("L33").Value = "Primera"
$ws.Range("M33").Value = 60
$ws.Range("D34").Value = 44937
$ws.Range("L34").Value = "Primera"
$ws.Range("M34").Value = 100
$ws.Range("O34").Value = 3000
$ws.Range("P34").Value = 2750
$ws.Range("S34").Value = 1375
$ws.Range("D35").Value = 44935
$ws.Range("M35").Value = 50
$ws.Range("N35").Value = 3000
$ws.Range("P35").Value = 3000
$ws.Range("S35").Value = 1500
$ws.Range("D36").Value = 44951
$ws.Range("D37").Value = 44972
$ws.Range("L37").Value = "Segunda"
$ws.Range("M37").Value = 30
$ws.Range("N37").Value = 2500
$ws.Range("O37").Value = 2500
$ws.Range("P37").Value = 2500
$ws.Range("S37").Value = 1250
$ws.Range("D38").Value = 44187
$ws.Range("M38").Value = 80
$ws.Range("N38").Value = 2800
$ws.Range("P38").Value = 2900
$ws.Range("R38").Value = "Provincia de Linares"
$ws.Range("S38").Value = 1450
$ws.Range("D39").Value = 44187
$ws.Range("L39").Value = "Primera"
$ws.Range("M39").Value = 65
$ws.Range("N39").Value = 1400
$ws.Range("O39").Value = 1500
$ws.Range("P39").Value = 1446
$ws.Range("Q39").Value = "`$/envase 1 kilo"
$ws.Range("S39").Value = 1446
$ws.Range("T39").Value = 1
$ws.Range("D40").Value = 44932
$ws.Range("M40").Value = 60
$ws.Range("N40").Value = 3000
$ws.Range("O40").Value = 3000
$ws.Range("P40").Value = 3000
$ws.Range("R40").Value = "Provincia de Diguillín"
$ws.Range("S40").Value = 1500
$ws.Range("D41").Value = 44952
$ws.Range("M41").Value = 30
$ws.Range("N41").Value = 3000
$ws.Range("O41").Value = 3000
$ws.Range("P41").Value = 3000
$ws.Range("S41").Value = 1500
$ws.Range("D42").Value = 44953
$ws.Range("M42").Value = 30
$ws.Range("N42").Value = 3000
$ws.Range("O42").Value = 3000
$ws.Range("P42").Value = 3000
$ws.Range("Q42").Value = "`$/bandeja 2 kilos"
$ws.Range("S42").Value = 1500
$ws.Range("T42").Value = 2
$ws.Range("D43").Value = 44979
$ws.Range("M43").Value = 30
$ws.Range("D44").Value = 44979
$ws.Range("D45").Value = 44594
$ws.Range("L45").Value = "Primera"
$ws.Range("M45").Value = 120
$ws.Range("O45").Value = 2800
$ws.Range("P45").Value = 2650
$ws.Range("R45").Value = "Provincia de Linares"
$ws.Range("S45").Value = 1325
